# ---------------------------------------------------------------------------
# Edit summary (per the XML diff / commit "Added test photos and created
# code for json file"):
#   - Swap the "Yash Jain" (Id 3) and "Logan Singh" (Id 4) rows, which were
#     rows 2 and 3 respectively (Logan Singh's full row -- including its
#     Photo hyperlink -- now comes first as row 2; Yash Jain becomes row 3).
#   - Fix the "Justin Tran" name in row 5: was stored with two trailing
#     spaces ("Justin Tran  "), now trimmed to "Justin Tran".
#   - Move the active-cell selection to C15.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 2 (Yash Jain / Id 3) and row 3 (Logan Singh / Id 4) ----------
# Columns: A=Id, B=Start time, C=Completion time, D=Email, E=Name,
#          F=HandleIG, G=Photo (hyperlink; blank for Yash Jain's row).

# Row 2 currently (before edit): Yash Jain / Id 3
$r2_id      = 3
$r2_start   = 45566.654479166667
$r2_end     = 45566.654583333337
$r2_email   = "yash.jain@surreyschools.ca"
$r2_name    = "Yash Jain"
$r2_handle  = "test1"

# Row 3 currently (before edit): Logan Singh / Id 4
$r3_id      = 4
$r3_start   = 45644.429293981499
$r3_end     = 45644.429826388892
$r3_email   = "logan.singh@surreyschools.ca"
$r3_name    = "Logan Singh"
$r3_handle  = "test2"
$r3_photo   = "https://sd36-my.sharepoint.com/personal/yash_jain_surreyschools_ca/Documents/Apps/Microsoft%20Forms/Grad%20Tag/Question/20241204_221947_Logan%20Singh.jpg"

# Write Logan Singh's data into row 2, Yash Jain's data into row 3.
$ws.Range("A2").Value = $r3_id
$ws.Range("B2").Value = $r3_start
$ws.Range("C2").Value = $r3_end
$ws.Range("D2").Value = $r3_email
$ws.Range("E2").Value = $r3_name
$ws.Range("F2").Value = $r3_handle
$ws.Range("G2").Value = $r3_photo

$ws.Range("A3").Value = $r2_id
$ws.Range("B3").Value = $r2_start
$ws.Range("C3").Value = $r2_end
$ws.Range("D3").Value = $r2_email
$ws.Range("E3").Value = $r2_name
$ws.Range("F3").Value = $r2_handle
$ws.Range("G3").Value = ""

# --- Rebuild hyperlinks to match the swapped rows --------------------------
# Deleting hyperlinks anywhere on the sheet clears the whole collection in
# this host, so remove them all once, then recreate the four links pointing
# at their (possibly new) cells, in the same order as the target file.
$ws.Range("D2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:yash.jain@surreyschools.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://sd36-my.sharepoint.com/personal/yash_jain_surreyschools_ca/Documents/Apps/Microsoft%20Forms/Grad%20Tag/Question/20241204_221947_Logan%20Singh.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://sd36-my.sharepoint.com/personal/yash_jain_surreyschools_ca/Documents/Apps/Microsoft%20Forms/Grad%20Tag/Question/IMG_5471_Harangad%20Sidhu.jpeg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "https://sd36-my.sharepoint.com/personal/yash_jain_surreyschools_ca/Documents/Apps/Microsoft%20Forms/Grad%20Tag/Question/image_Justin%20Tran.jpg") | Out-Null

# --- Fix trailing whitespace in "Justin Tran" (row 5, Name column) ---------
$ws.Range("E5").Value = "Justin Tran"

# --- Update the sheet's active-cell selection -------------------------------
$ws.Range("C15").Select()
